# Populate the "Customer" test-data sheet (sheet1 / rId1) with header row +
# one sample data row, matching the "data sheet modified and added to
# tests" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")

# Make it the active sheet/tab (matches tabSelected moving from "Reason" to
# "Customer" in the target workbook).
[void]$ws.Select()

# --- Header row (bold) ---------------------------------------------------
$ws.Range("A1").Value = "customer_name"
$ws.Range("B1").Value = "customer_code"
$ws.Range("C1").Value = "customer_period"
$ws.Range("D1").Value = "customer_drop"
$ws.Range("A1:D1").Font.Bold = $true

# --- Data row --------------------------------------------------------------
$ws.Range("A2").NumberFormat = "0"
$ws.Range("A2").Value = "Akash"
$ws.Range("B2").Value = 2907657
$ws.Range("C2").Value = 14
$ws.Range("D2").NumberFormat = "0"
$ws.Range("D2").Value = "AOB"

# --- Column widths (approximate the authored layout; the host's column
# -width unit only quantizes to ~1/6-character steps, so these are the
# closest achievable values to the authored 17.82/15.18/15.09/16.91) --------
$ws.Columns.Item(1).ColumnWidth = 16.9214
$ws.Columns.Item(2).ColumnWidth = 14.2547
$ws.Columns.Item(3).ColumnWidth = 14.2548
$ws.Columns.Item(4).ColumnWidth = 15.9262

# --- Leave the last selection where the original author left it ------------
[void]$ws.Range("E6").Select()
